# edit.ps1 -- apply the Makale Yonetim Sistemi Proje Raporu revision
# Described by the commit "Add files via upload":
#  - Split the single Frontend/.../Diger Araclar paragraph so that the
#    former "Backend:" bullet list is replaced by the new Flask-based
#    backend stack (Flask, Flask-RESTful, Flask-JWT-Extended, SQLAlchemy),
#    each new line living in its own paragraph with an Arial rFonts rPr,
#    plus w:proofErr spell-check markers around the borrowed English terms.
#  - Wrap "MySQL :" in gramStart/gramEnd proofErr marks and split the
#    "- " prefix into its own run.
#  - Drop the two stray <w:lastRenderedPageBreak/> runs that no longer
#    apply once the content reflows.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Replace the "Frontend: ... Diger Araclar: ..." paragraph with the
#    new multi-paragraph Backend stack listing.
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("Frontend:", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$target = $rng.Paragraphs(1).Range

$newBlockXml = @'
<?xml version="1.0" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p><w:pPr><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr></w:pPr><w:r><w:t>Frontend:</w:t></w:r><w:r><w:br/><w:t>- React: Kullanıcı arayüzü geliştirme.</w:t></w:r><w:r><w:br/><w:t>- React Router: Sayfa geçişleri ve yönlendirme.</w:t></w:r><w:r><w:br/><w:t>- TypeScript: Daha güvenilir bir kod yapısı için.</w:t></w:r><w:r><w:br/><w:t>- Bootstrap: Duyarlı ve modern tasarım için CSS framework.</w:t></w:r><w:r><w:br/></w:r><w:r><w:br/></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr><w:t>Backend</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr><w:t>:</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr><w:t xml:space="preserve">- </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr><w:t>Flask</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr><w:t xml:space="preserve">: </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr><w:t>RESTful</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr><w:t xml:space="preserve"> API geliştirme ve sunucu tarafı işlemleri yönetme.</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr><w:t xml:space="preserve">- </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr><w:t>Flask-RESTful</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr><w:t xml:space="preserve">: API </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr><w:t>endpoint’lerinin</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr><w:t xml:space="preserve"> yapılandırılması.</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr><w:t xml:space="preserve">- </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr><w:t>Flask</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr><w:t>-JWT-</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr><w:t>Extended</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr><w:t>: Kullanıcı kimlik doğrulama ve oturum yönetimi.</w:t></w:r></w:p><w:p><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr><w:t xml:space="preserve">- </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr><w:t>SQLAlchemy</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr><w:t xml:space="preserve">: </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr><w:t>Veritabanı</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr><w:t xml:space="preserve"> modelleme ve ORM işlemleri.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr><w:br/></w:r><w:r><w:br/><w:t>Veritabanı:</w:t></w:r><w:r><w:br/><w:t xml:space="preserve">- </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>MySQL :</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> Makale ve kullanıcı verilerinin depolanması.</w:t></w:r><w:r><w:br/></w:r><w:r><w:br/><w:t>Diğer Araçlar:</w:t></w:r><w:r><w:br/><w:t>- Axios: Frontend ile backend arasında veri alışverişi için HTTP istekleri.</w:t></w:r></w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@

$target.InsertXML($newBlockXml)

# ---------------------------------------------------------------------
# 2) Remove <w:lastRenderedPageBreak/> before "3. Proje Yapisi".
# ---------------------------------------------------------------------
$rng2 = $d.Content
$rng2.Find.Execute("3. Proje Yap", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$heading = $rng2.Paragraphs(1).Range

$headingXml = @'
<?xml version="1.0" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p><w:pPr><w:rPr><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="4C94D8" w:themeColor="text2" w:themeTint="80"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t>3. Proje Yapısı</w:t></w:r></w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@

$heading.InsertXML($headingXml)

# ---------------------------------------------------------------------
# 3) Remove <w:lastRenderedPageBreak/> from the blank break/break
#    paragraph that precedes the "Frontend Klasor Yapisi" heading.
# ---------------------------------------------------------------------
$rng3 = $d.Content
$rng3.Find.Execute("Frontend Klas", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$breakPara = $rng3.Paragraphs(1).Previous().Previous()
$breakRng = $breakPara.Range

$breakXml = @'
<?xml version="1.0" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p><w:pPr><w:pStyle w:val="Balk3"/><w:rPr><w:color w:val="FF0000"/></w:rPr></w:pPr><w:r><w:br/></w:r><w:r><w:br/></w:r></w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@

$breakRng.InsertXML($breakXml)

Write-Output "Edit complete."
